# Update the division problems shown in the worksheet table.
#
# Cells are addressed by (row, column) rather than by searching for the
# old text, because:
#   1) Several of the old values are duplicated across the table
#      (e.g. "65÷5=" appears twice), so a plain Find/Replace could not
#      tell the occurrences apart.
#   2) Assigning directly to a table cell's Range.Text reliably scopes
#      the edit to just that cell while preserving the existing run
#      formatting (font/size) defined on that cell's paragraph.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; New="17÷9="},
    @{Row=1;  Col=2; New="50÷2="},
    @{Row=1;  Col=3; New="70÷2="},
    @{Row=1;  Col=4; New="14÷9="},
    @{Row=1;  Col=5; New="66÷3="},

    @{Row=5;  Col=1; New="63÷6="},
    @{Row=5;  Col=2; New="20÷8="},
    @{Row=5;  Col=3; New="28÷7="},
    @{Row=5;  Col=4; New="25÷7="},
    @{Row=5;  Col=5; New="22÷7="},

    @{Row=9;  Col=1; New="70÷9="},
    @{Row=9;  Col=2; New="76÷5="},
    @{Row=9;  Col=3; New="37÷6="},
    @{Row=9;  Col=4; New="83÷3="},
    @{Row=9;  Col=5; New="40÷5="},

    @{Row=13; Col=1; New="55÷5="},
    @{Row=13; Col=2; New="19÷5="},
    @{Row=13; Col=3; New="77÷3="},
    @{Row=13; Col=4; New="73÷2="},
    @{Row=13; Col=5; New="56÷4="},

    @{Row=17; Col=1; New="64÷4="},
    @{Row=17; Col=2; New="33÷2="},
    @{Row=17; Col=3; New="40÷7="},
    @{Row=17; Col=4; New="38÷5="},
    @{Row=17; Col=5; New="45÷7="}
)

foreach ($rep in $replacements) {
    $cell = $t.Cell($rep.Row, $rep.Col)
    $cell.Range.Text = $rep.New
}
